$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = [DateTime]"2023-06-30"
$ws.Range("C5").Value = "Activision Blizzard"
$ws.Range("D5").Value = "ATVI"
$ws.Range("D5").WrapText = $true
$ws.Range("E5").Value = 3861.24
$ws.Range("F5").Value = 83.94

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:I11")) | Out-Null

$ws.Range("J5").Select() | Out-Null
